$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 213.5
$ws.Range("I6").Value = 156.2
$ws.Range("K6").Value = 468.6
$ws.Range("M6").Value = -356.6
$ws.Range("H40").Value = 1833
$ws.Range("I40").Value = 1620.7858
$ws.Range("J40").Value = 2080.5833
$ws.Range("K40").Value = 1620.7858
$ws.Range("L40").Value = 2080.5833
$ws.Range("M40").Value = -1445.7858
$ws.Range("N40").Value = -2430.5833
$ws.Range("H132").Value = 2743.2222
$ws.Range("I132").Value = 2172.2654
$ws.Range("J132").Value = 4741.5713
$ws.Range("K132").Value = 6516.796200000001
$ws.Range("L132").Value = 14224.7139
$ws.Range("M132").Value = -3986.796200000001
$ws.Range("N132").Value = -19284.7139
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H135").Value = 311.33334
$ws.Range("I135").Value = 311.33334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2802.00006
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -267.0000600000003
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 358941.47
$ws.Range("I61").Value = 1915.25
$ws.Range("J61").Value = 1251507
$ws.Range("K61").Value = 1915.25
$ws.Range("L61").Value = 1251507
$ws.Range("M61").Value = -1703.25
$ws.Range("N61").Value = -1251931
$ws.Range("H63").Value = 1754
$ws.Range("I63").Value = 1899.3334
$ws.Range("J63").Value = 1463.3334
$ws.Range("K63").Value = 1899.3334
$ws.Range("L63").Value = 1463.3334
$ws.Range("M63").Value = -1213.3334
$ws.Range("N63").Value = -2835.3334
$ws.Range("H66").Value = 1754
$ws.Range("I66").Value = 1899.3334
$ws.Range("J66").Value = 1463.3334
$ws.Range("K66").Value = 9496.666999999999
$ws.Range("L66").Value = 7316.666999999999
$ws.Range("M66").Value = -6064.666999999999
$ws.Range("N66").Value = -14180.667
$ws.Range("H132").Value = 4213.8945
$ws.Range("I132").Value = 1067.4286
$ws.Range("J132").Value = 13024
$ws.Range("K132").Value = 3202.2858
$ws.Range("L132").Value = 39072
$ws.Range("M132").Value = -672.2857999999997
$ws.Range("N132").Value = -44132
$ws.Range("H134").Value = 58888
$ws.Range("J134").Value = 58888
$ws.Range("L134").Value = 58888
$ws.Range("N134").Value = -69028
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 358941.47
$ws.Range("I136").Value = 1915.25
$ws.Range("J136").Value = 1251507
$ws.Range("K136").Value = 5745.75
$ws.Range("L136").Value = 3754521
$ws.Range("M136").Value = -3195.75
$ws.Range("N136").Value = -3759621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1576.0454
$ws.Range("I107").Value = 1155.4706
$ws.Range("J107").Value = 3006
$ws.Range("K107").Value = 1155.4706
$ws.Range("L107").Value = 3006
$ws.Range("M107").Value = 764.5293999999999
$ws.Range("N107").Value = -6846
$ws.Range("H132").Value = 70780
$ws.Range("J132").Value = 70780
$ws.Range("L132").Value = 70780
$ws.Range("N132").Value = -80900
$ws.Range("H134").Value = 1313.7188
$ws.Range("I134").Value = 985.88
$ws.Range("J134").Value = 2484.5715
$ws.Range("K134").Value = 2957.64
$ws.Range("L134").Value = 7453.7145
$ws.Range("M134").Value = -422.6399999999999
$ws.Range("N134").Value = -12523.7145
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4181.0586
$ws.Range("I62").Value = 4678
$ws.Range("J62").Value = 2800.6667
$ws.Range("K62").Value = 4678
$ws.Range("L62").Value = 2800.6667
$ws.Range("M62").Value = -4054
$ws.Range("N62").Value = -4048.6667
$ws.Range("H65").Value = 4181.0586
$ws.Range("I65").Value = 4678
$ws.Range("J65").Value = 2800.6667
$ws.Range("K65").Value = 23390
$ws.Range("L65").Value = 14003.3335
$ws.Range("M65").Value = -20270
$ws.Range("N65").Value = -20243.3335
$ws.Range("H134").Value = 2190.7576
$ws.Range("I134").Value = 2287.24
$ws.Range("J134").Value = 1889.25
$ws.Range("K134").Value = 6861.719999999999
$ws.Range("L134").Value = 5667.75
$ws.Range("M134").Value = -4326.719999999999
$ws.Range("N134").Value = -10737.75
$ws.Range("H135").Value = 50668.57
$ws.Range("J135").Value = 50668.57
$ws.Range("L135").Value = 50668.57
$ws.Range("N135").Value = -60808.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1000.3333
$ws.Range("I16").Value = 1000.3333
$ws.Range("K16").Value = 3000.9999
$ws.Range("M16").Value = -2827.9999
$ws.Range("H92").Value = 435.2
$ws.Range("I92").Value = 312
$ws.Range("J92").Value = 517.3333
$ws.Range("K92").Value = 936
$ws.Range("L92").Value = 1551.9999
$ws.Range("M92").Value = 312
$ws.Range("N92").Value = -4047.9999
$ws.Range("H97").Value = 429
$ws.Range("I97").Value = 335.33334
$ws.Range("J97").Value = 452.41666
$ws.Range("K97").Value = 1006.00002
$ws.Range("L97").Value = 1357.24998
$ws.Range("M97").Value = -510.0000200000001
$ws.Range("N97").Value = -2349.24998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31080
$ws.Range("H113").Value = 2001182.2
$ws.Range("I113").Value = 2500977.8
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2500977.8
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -2498807.8
$ws.Range("N113").Value = -6340
$ws.Range("H114").Value = 42333.332
$ws.Range("J114").Value = 42333.332
$ws.Range("L114").Value = 42333.332
$ws.Range("N114").Value = -51011.332
$ws.Range("H122").Value = 2480.8
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 4302.6665
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 12907.9995
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -17807.9995
$ws.Range("H132").Value = 38215.5
$ws.Range("I132").Value = 60397.234
$ws.Range("J132").Value = 3934.6365
$ws.Range("K132").Value = 181191.702
$ws.Range("L132").Value = 11803.9095
$ws.Range("M132").Value = -178661.702
$ws.Range("N132").Value = -16863.9095
$ws.Range("H133").Value = 38400
$ws.Range("J133").Value = 38400
$ws.Range("L133").Value = 38400
$ws.Range("N133").Value = -48520
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H132").Value = 22095.79
$ws.Range("I132").Value = 4690.3
$ws.Range("J132").Value = 41435.223
$ws.Range("K132").Value = 14070.9
$ws.Range("L132").Value = 124305.669
$ws.Range("M132").Value = -11540.9
$ws.Range("N132").Value = -129365.669
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 42000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -52140
$ws.Range("H137").Value = 30290
$ws.Range("I137").Value = 30290
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 30290
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -25190
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 1100
$ws.Range("I9").Value = 1100
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1100
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -960
$ws.Range("N9").ClearContents()
$ws.Range("H132").Value = 1706.6666
$ws.Range("I132").Value = 1545.9231
$ws.Range("J132").Value = 1855.9286
$ws.Range("K132").Value = 4637.7693
$ws.Range("L132").Value = 5567.7858
$ws.Range("M132").Value = -2107.7693
$ws.Range("N132").Value = -10627.7858
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
